$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '27.083.24'
$ws.Cells.Item(2, 5).Value = '  +0.02%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.832.22'
$ws.Cells.Item(3, 5).Value = '  +0.41%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '1.007'
$ws.Cells.Item(4, 5).Value = '  -0.13%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '311.56'
$ws.Cells.Item(5, 5).Value = '  -0.24%  '

# Row 6
$ws.Cells.Item(6, 5).Value = '  +0.07%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '0.4656'
$ws.Cells.Item(7, 5).Value = '  -0.51%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.3709'
$ws.Cells.Item(8, 5).Value = '  +1.63%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.07374'
$ws.Cells.Item(9, 5).Value = '  -0.05%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '0.8746'
$ws.Cells.Item(10, 5).Value = '  -0.42%  '

# Row 11
$ws.Cells.Item(11, 2).Value = 'Solana'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(11, 4).Value = '19.99'
$ws.Cells.Item(11, 5).Value = '  -1.19%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).Value = '0.07852'
$ws.Cells.Item(12, 5).Value = '  +5.51%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'Chainlink'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(13, 4).Value = '6.629'
$ws.Cells.Item(13, 5).Value = '  +1.68%  '

# Row 14
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).Value = '5.354'
$ws.Cells.Item(14, 5).Value = '  -0.26%  '

# Row 15
$ws.Cells.Item(15, 2).Value = 'Litecoin'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(15, 4).Value = '91.91'
$ws.Cells.Item(15, 5).Value = '  -1.04%  '

# Row 16
$ws.Cells.Item(16, 2).Value = 'WrappedEther'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(16, 4).Value = '1.707.50'
$ws.Cells.Item(16, 5).Value = '  -9.96%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '1.008'
$ws.Cells.Item(17, 5).Value = '  +0.14%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '0.000008862'
$ws.Cells.Item(18, 5).Value = '  +1.71%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '1.009'
$ws.Cells.Item(19, 5).Value = '  +0.22%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'WrappedBTC'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(20, 4).Value = '27.459.93'
$ws.Cells.Item(20, 5).Value = '  -0.37%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'Avalanche'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(21, 4).Value = '14.65'
$ws.Cells.Item(21, 5).Value = '  +0.41%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '5.139'
$ws.Cells.Item(22, 5).Value = '  -1.75%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '10.58'
$ws.Cells.Item(23, 5).Value = '  -0.23%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '1.965.96'
$ws.Cells.Item(24, 5).Value = '  -5.52%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '152.36'
$ws.Cells.Item(25, 5).Value = '  +0.72%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '1.827'
$ws.Cells.Item(26, 5).Value = '  -2.85%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '18.31'
$ws.Cells.Item(27, 5).Value = '  -0.92%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '2.097'
$ws.Cells.Item(28, 5).Value = '  -1.64%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.01%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '115.27'
$ws.Cells.Item(30, 5).Value = '  -0.86%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '0.08885'
$ws.Cells.Item(31, 5).Value = '  -0.32%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '2.964'
$ws.Cells.Item(32, 5).Value = '  +0.83%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(33, 4).Value = '0.7289'
$ws.Cells.Item(33, 5).Value = '  -1.95%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'Filecoin'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(34, 4).Value = '4.442'
$ws.Cells.Item(34, 5).Value = '  -1.36%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.139'
$ws.Cells.Item(35, 5).Value = '  -1.94%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '2.503'
$ws.Cells.Item(36, 5).Value = '  -0.84%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(37, 4).Value = '1.074'
$ws.Cells.Item(37, 5).Value = '  -1.36%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'VeChain'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(38, 4).Value = '0.01951'
$ws.Cells.Item(38, 5).Value = '  +1.03%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.05224'
$ws.Cells.Item(39, 5).Value = '  -1.26%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '2.933'
$ws.Cells.Item(40, 5).Value = '  -0.01%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '7.186'
$ws.Cells.Item(41, 5).Value = '  -2.16%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '0.5203'
$ws.Cells.Item(42, 5).Value = '  -0.82%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.8844'
$ws.Cells.Item(43, 5).Value = '  -12.15%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1630'
$ws.Cells.Item(44, 5).Value = '  -0.47%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.240'
$ws.Cells.Item(45, 5).Value = '  -1.40%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '0.4827'
$ws.Cells.Item(46, 5).Value = '  -1.28%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +0.15%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '10.14'
$ws.Cells.Item(48, 5).Value = '  -2.76%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '102.59'
$ws.Cells.Item(49, 5).Value = '  -1.65%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '1.630'
$ws.Cells.Item(50, 5).Value = '  -1.13%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.06216'
$ws.Cells.Item(51, 5).Value = '  -0.83%  '
